# televisies.xlsx — "Wijzigingen en aanvullingen prognoses"
#
# The source sheet was resaved by a newer Excel build and tidied up a bit:
#   * the worksheet was renamed from the generic "Blad1" to "data"
#   * the sheet's cursor/selection cache was cleared (no more <selection>
#     left over from editing cell B1)
#   * a page setup (A4 / portrait) was stamped onto the sheet for printing
#
# Apply the parts of that which are real, user-visible workbook state and
# reachable from the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab: "Blad1" -> "data"
$ws.Name = "data"

# 2. Reset the active selection back to the top-left cell so the sheet
#    doesn't keep carrying around the old "B1 selected" cursor state.
$null = $ws.Range("A1").Select()

# 3. Give the sheet an explicit print setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = [Microsoft.Office.Interop.Excel.XlPaperSize]::xlPaperA4
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait
